$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its original text formatting so that
# numeric-looking values (e.g. "1.003", "6.000") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "27.955.79"
$ws.Cells.Item(2, 5).Value = "  +0.07%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.858.59"
$ws.Cells.Item(3, 5).Value = "  -0.70%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.004"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "311.44"
$ws.Cells.Item(5, 5).Value = "  -0.29%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "1.003"
$ws.Cells.Item(6, 5).Value = "  +0.01%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5126"
$ws.Cells.Item(7, 5).Value = "  +3.00%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.24%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.08271"
$ws.Cells.Item(9, 5).Value = "  -8.45%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(10, 4).Value = "1.108"
$ws.Cells.Item(10, 5).Value = "  -0.73%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "OKB"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(11, 4).Value = "41.58"
$ws.Cells.Item(11, 5).Value = "  +0.14%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "6.185"
$ws.Cells.Item(12, 5).Value = "  -1.66%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.863.54"
$ws.Cells.Item(13, 5).Value = "  -1.23%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "20.45"
$ws.Cells.Item(14, 5).Value = "  -0.96%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "7.192"
$ws.Cells.Item(15, 5).Value = "  -0.34%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -0.05%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.00001092"
$ws.Cells.Item(17, 5).Value = "  -0.70%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "90.34"
$ws.Cells.Item(18, 5).Value = "  -0.53%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.06600"
$ws.Cells.Item(19, 5).Value = "  -0.64%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "17.74"
$ws.Cells.Item(20, 5).Value = "  -0.69%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.01%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "6.000"
$ws.Cells.Item(22, 5).Value = "  -1.96%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "27.987.09"
$ws.Cells.Item(23, 5).Value = "  -0.02%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "11.01"
$ws.Cells.Item(24, 5).Value = "  -3.10%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "2.211"
$ws.Cells.Item(25, 5).Value = "  -3.25%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.578"
$ws.Cells.Item(26, 5).Value = "  +2.43%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "2.074.75"
$ws.Cells.Item(27, 5).Value = "  -1.36%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "157.01"
$ws.Cells.Item(28, 5).Value = "  -0.41%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "20.36"
$ws.Cells.Item(29, 5).Value = "  -1.96%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "124.45"
$ws.Cells.Item(30, 5).Value = "  -1.65%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.1061"
$ws.Cells.Item(31, 5).Value = "  +0.99%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "1.035"
$ws.Cells.Item(32, 5).Value = "  -1.71%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "5.600"
$ws.Cells.Item(33, 5).Value = "  +0.62%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "3.606"
$ws.Cells.Item(34, 5).Value = "  +0.41%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "9.523"
$ws.Cells.Item(35, 5).Value = "  +2.45%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "0.06536"
$ws.Cells.Item(36, 5).Value = "  -0.11%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).Value = "0.02419"
$ws.Cells.Item(37, 5).Value = "  +0.89%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.2173"
$ws.Cells.Item(38, 5).Value = "  -0.28%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "1.207"
$ws.Cells.Item(39, 5).Value = "  +1.13%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "0.6418"
$ws.Cells.Item(40, 5).Value = "  +0.68%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -2.90%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -2.41%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "4.871"
$ws.Cells.Item(43, 5).Value = "  -0.47%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "0.6089"
$ws.Cells.Item(44, 5).Value = "  +1.35%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "13.11"
$ws.Cells.Item(45, 5).Value = "  -1.19%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.73%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "3.649"
$ws.Cells.Item(47, 5).Value = "  -0.98%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "1.972"
$ws.Cells.Item(48, 5).Value = "  -0.16%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "1.203"
$ws.Cells.Item(49, 5).Value = "  -0.70%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "120.66"
$ws.Cells.Item(50, 5).Value = "  +0.42%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "79.23"
$ws.Cells.Item(51, 5).Value = "  +0.91%  "
